# Applies the "Updated notebook, reran simulation" edit:
#  - Inserts two new categories ("Holden", "Rizzie Spiral") right after "Spiral5"
#    (pushing all subsequent category rows down by 2 rows).
#  - Renames "Thomas Hex" -> "Matthies Hex".
#  - Refreshes the simulated numeric results (columns C:W) for every category,
#    including brand-new rows for the two newly added categories.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for the two new categories by inserting two blank rows right
#    after the "Spiral5" row (row 3), i.e. before the old row 4.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Resize(2).Insert()

# The insert leaves the new A4:A5 cells without the bold/bordered "index"
# style used throughout column A (e.g. A6) - copy that formatting over.
$ws.Range("A6").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# 2. Renumber column A (the 0-based category index) for every row that was
#    pushed down by the insert: old A6..A31 become A8..A33 -> i.e. add 2 to
#    every existing index from row 6 downward.
# ---------------------------------------------------------------------------
$lastRow = 31
for ($r = 6; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cur = $cell.Value()
    if ($null -ne $cur) {
        $cell.Value = $cur + 2
    }
}

# ---------------------------------------------------------------------------
# 3. Populate the two new rows (row 4 = "Holden", row 5 = "Rizzie Spiral").
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"

$row4 = @{
    "C" = 0.9984268265051536
    "D" = 1.000604876613467
    "E" = 1.000898099160029
    "F" = 0.9994801070681794
    "G" = 1.001613000642363
    "H" = 1.001613000642363
    "I" = 1.001613000642363
    "J" = 0.9975805017029581
    "K" = 1.000604876613467
    "L" = 0.9995160998379862
    "M" = 0.9975805017029581
    "N" = 1.001613000642363
    "O" = 1.000604876613467
    "P" = 0.9990926891582128
    "Q" = 1.000042491840823
    "R" = 0.9999327929862628
    "S" = 0.999221828461535
    "T" = 0.9999327929862627
    "U" = 0.9998196215067419
    "V" = 1.000178297333866
    "W" = 0.9998405485179505
}

$row5 = @{
    "C" = 0.9888219340638643
    "D" = 1.004297895656907
    "E" = 1.00638137480451
    "F" = 0.9963059408538191
    "G" = 1.011461061223646
    "H" = 1.011461061223646
    "I" = 1.011461061223646
    "J" = 0.9828084162066943
    "K" = 1.004297895656907
    "L" = 0.9965616828357537
    "M" = 0.9828084162066943
    "N" = 1.011461061223646
    "O" = 1.004297895656907
    "P" = 0.9935531559318007
    "Q" = 1.000301918255363
    "R" = 0.9995224576957492
    "S" = 0.9944707509058069
    "T" = 0.9995224576957492
    "U" = 0.9987183284852668
    "V" = 1.001266875032943
    "W" = 0.9988670251627628
}

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")
foreach ($c in $cols) {
    $ws.Range($c + "4").Value = $row4[$c]
    $ws.Range($c + "5").Value = $row5[$c]
}

# ---------------------------------------------------------------------------
# 4. Rename the "Thomas Hex" category (now on row 11 after the insert) to
#    "Matthies Hex".
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "Matthies Hex"

# The used range / dimension (A1:W31) updates automatically as rows are
# inserted and populated above.
$wb.Save()
